$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 715.6286
$ws.Range("J33").Value = 1570.7142
$ws.Range("L33").Value = 1570.7142
$ws.Range("N33").Value = -2028.7142

$ws.Range("H86").Value = 300030000
$ws.Range("I86").Value = 350034340
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 350034340
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -350033217
$ws.Range("N86").Value = -6246

$ws.Range("H89").Value = 300030000
$ws.Range("I89").Value = 350034340
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 1750171700
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -1750166084
$ws.Range("N89").Value = -31232

$ws.Range("H98").Value = 1112.5
$ws.Range("I98").Value = 300
$ws.Range("K98").Value = 300
$ws.Range("M98").Value = 1198

$ws.Range("H111").Value = 7866.375
$ws.Range("I111").Value = 9676.333000000001
$ws.Range("J111").Value = 6780.4
$ws.Range("K111").Value = 29028.999
$ws.Range("L111").Value = 20341.2
$ws.Range("M111").Value = -25961.999
$ws.Range("N111").Value = -26475.2

$ws.Range("H112").Value = 4690.163
$ws.Range("J112").Value = 4823.787
$ws.Range("L112").Value = 14471.361
$ws.Range("N112").Value = -16687.361

$ws.Range("H122").Value = 1112.5
$ws.Range("I122").Value = 300
$ws.Range("K122").Value = 900
$ws.Range("M122").Value = 1550

$ws.Range("H138").Value = 2125.3333
$ws.Range("I138").Value = 2184.1
$ws.Range("J138").Value = 2102.7307
$ws.Range("K138").Value = 6552.299999999999
$ws.Range("L138").Value = 6308.1921
$ws.Range("M138").Value = -1412.299999999999
$ws.Range("N138").Value = -16588.1921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 363683.06
$ws.Range("I32").Value = 406453.53
$ws.Range("J32").Value = 21519.3
$ws.Range("K32").Value = 406453.53
$ws.Range("L32").Value = 21519.3
$ws.Range("M32").Value = -406166.53
$ws.Range("N32").Value = -22093.3

$ws.Range("H132").Value = 1813.9348
$ws.Range("I132").Value = 1390.3513
$ws.Range("J132").Value = 3555.3333
$ws.Range("K132").Value = 4171.0539
$ws.Range("L132").Value = 10665.9999
$ws.Range("M132").Value = -1641.0539
$ws.Range("N132").Value = -15725.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 111111110
$ws.Range("I7").Value = 111111110
$ws.Range("K7").Value = 111111110
$ws.Range("M7").Value = -111110997

$ws.Range("H38").Value = 10000
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10832

$ws.Range("H97").Value = 27499
$ws.Range("J97").Value = 34998
$ws.Range("L97").Value = 34998
$ws.Range("N97").Value = -36980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 29666.889
$ws.Range("J4").Value = 29666.889
$ws.Range("L4").Value = 29666.889
$ws.Range("N4").Value = -29890.889

$ws.Range("H16").Value = 2420
$ws.Range("I16").Value = 1830.8
$ws.Range("J16").Value = 3156.5
$ws.Range("K16").Value = 1830.8
$ws.Range("L16").Value = 3156.5
$ws.Range("M16").Value = -1543.8
$ws.Range("N16").Value = -3730.5

$ws.Range("H31").Value = 8117.758
$ws.Range("I31").Value = 1740.1111
$ws.Range("K31").Value = 1740.1111
$ws.Range("M31").Value = -1445.1111

$ws.Range("H34").Value = 8117.758
$ws.Range("I34").Value = 1740.1111
$ws.Range("K34").Value = 1740.1111
$ws.Range("M34").Value = -1538.1111

$ws.Range("H113").Value = 2420
$ws.Range("I113").Value = 1830.8
$ws.Range("J113").Value = 3156.5
$ws.Range("K113").Value = 1830.8
$ws.Range("L113").Value = 3156.5
$ws.Range("M113").Value = 339.2
$ws.Range("N113").Value = -7496.5

$ws.Range("H134").Value = 6072.5
$ws.Range("I134").Value = 3357.5
$ws.Range("K134").Value = 10072.5
$ws.Range("M134").Value = -7537.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 876.09375
$ws.Range("I131").Value = 275
$ws.Range("J131").Value = 1076.4584
$ws.Range("K131").Value = 825
$ws.Range("L131").Value = 3229.3752
$ws.Range("M131").Value = 4215
$ws.Range("N131").Value = -13309.3752

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 501000
$ws.Range("I113").Value = 1000000
$ws.Range("K113").Value = 1000000
$ws.Range("M113").Value = -997830

$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -44920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8182.56
$ws.Range("I22").Value = 494.7
$ws.Range("K22").Value = 494.7
$ws.Range("M22").Value = -199.7

$ws.Range("H27").Value = 8182.56
$ws.Range("I27").Value = 494.7
$ws.Range("K27").Value = 494.7
$ws.Range("M27").Value = -387.7

$ws.Range("H55").Value = 581.8421
$ws.Range("I55").Value = 278.44446
$ws.Range("J55").Value = 854.9
$ws.Range("K55").Value = 278.44446
$ws.Range("L55").Value = 854.9
$ws.Range("M55").Value = -105.44446
$ws.Range("N55").Value = -1200.9

$ws.Range("H61").Value = 2421.7273
$ws.Range("J61").Value = 3680
$ws.Range("L61").Value = 3680
$ws.Range("N61").Value = -4084

$ws.Range("H82").Value = 2212.875
$ws.Range("I82").Value = 1540.6
$ws.Range("J82").Value = 3333.3333
$ws.Range("K82").Value = 1540.6
$ws.Range("L82").Value = 3333.3333
$ws.Range("M82").Value = -1179.6
$ws.Range("N82").Value = -4055.3333

$ws.Range("H85").Value = 2212.875
$ws.Range("I85").Value = 1540.6
$ws.Range("J85").Value = 3333.3333
$ws.Range("K85").Value = 1540.6
$ws.Range("L85").Value = 3333.3333
$ws.Range("M85").Value = -292.5999999999999
$ws.Range("N85").Value = -5829.3333

$ws.Range("H93").Value = 15322.5
$ws.Range("I93").Value = 22116
$ws.Range("K93").Value = 22116
$ws.Range("M93").Value = -20868

$ws.Range("H100").Value = 2370.818
$ws.Range("I100").Value = 2012.8572
$ws.Range("J100").Value = 2997.25
$ws.Range("K100").Value = 2012.8572
$ws.Range("L100").Value = 2997.25
$ws.Range("M100").Value = -1471.8572
$ws.Range("N100").Value = -4079.25

$ws.Range("H113").Value = 2421.7273
$ws.Range("J113").Value = 3680
$ws.Range("L113").Value = 3680
$ws.Range("N113").Value = -8020

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5139.1665
$ws.Range("I81").Value = 6434.2
$ws.Range("J81").Value = 4214.143
$ws.Range("K81").Value = 12868.4
$ws.Range("L81").Value = 8428.286
$ws.Range("M81").Value = -11807.4
$ws.Range("N81").Value = -10550.286

$ws.Range("H84").Value = 5139.1665
$ws.Range("I84").Value = 6434.2
$ws.Range("J84").Value = 4214.143
$ws.Range("K84").Value = 64342
$ws.Range("L84").Value = 42141.43
$ws.Range("M84").Value = -59038
$ws.Range("N84").Value = -52749.43

$ws.Range("H107").Value = 535.6
$ws.Range("I107").Value = 486.14285
$ws.Range("K107").Value = 1458.42855
$ws.Range("M107").Value = 461.5714499999999

$ws.Range("H132").Value = 3705461.2
$ws.Range("I132").Value = 1487.7646
$ws.Range("J132").Value = 5954302.5
$ws.Range("K132").Value = 4463.293799999999
$ws.Range("L132").Value = 17862907.5
$ws.Range("M132").Value = -1933.293799999999
$ws.Range("N132").Value = -17867967.5

